$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance Audit")
$ws.Activate()

# Update "Education capital" -> "Education + capital" (row 22)
$ws.Range("A22").Value = "Education + capital"

# Insert a new row below it for "Year 7 + transition" (new row 23)
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "Year 7 + transition"

# Update "Debt sustainability" -> "Debt + sustainability" (now row 35 after insertion above)
$ws.Range("A35").Value = "Debt + sustainability"

# Update "Courts backlog" -> "Courts + backlog" (now row 43 after insertion above)
$ws.Range("A43").Value = "Courts + backlog"

# Insert two new rows below it for the disability funding terms (new rows 44 and 45)
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(45).Insert()
$ws.Range("A44").Value = "Disability funding + school "
$ws.Range("A45").Value = "Disability funding + education"

# Update the active selection to match the author's cursor position after editing
$ws.Range("C8").Select() | Out-Null
